$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '25.999.95'
$ws.Cells.Item(2, 5).Value = '  +0.16%  '
$ws.Cells.Item(3, 4).Value = '1.641.36'
$ws.Cells.Item(3, 5).Value = '  -0.37%  '
$ws.Cells.Item(4, 4).Value = "'" + '1.002'
$ws.Cells.Item(4, 5).Value = '  -0.67%  '
$ws.Cells.Item(5, 4).Value = "'" + '215.13'
$ws.Cells.Item(5, 5).Value = '  -0.46%  '
$ws.Cells.Item(6, 4).Value = "'" + '0.5058'
$ws.Cells.Item(6, 5).Value = '  -0.90%  '
$ws.Cells.Item(7, 4).Value = "'" + '1.001'
$ws.Cells.Item(7, 5).Value = '  -0.50%  '
$ws.Cells.Item(8, 4).Value = "'" + '0.2581'
$ws.Cells.Item(8, 5).Value = '  +0.08%  '
$ws.Cells.Item(9, 4).Value = "'" + '0.06359'
$ws.Cells.Item(9, 5).Value = '  -1.03%  '
$ws.Cells.Item(10, 4).Value = "'" + '19.85'
$ws.Cells.Item(10, 5).Value = '  +0.87%  '
$ws.Cells.Item(11, 4).Value = "'" + '0.07754'
$ws.Cells.Item(11, 5).Value = '  -0.53%  '
$ws.Cells.Item(12, 4).Value = "'" + '4.275'
$ws.Cells.Item(12, 5).Value = '  -1.11%  '
$ws.Cells.Item(13, 4).Value = '1.636.33'
$ws.Cells.Item(13, 5).Value = '  -0.69%  '
$ws.Cells.Item(14, 4).Value = "'" + '0.5482'
$ws.Cells.Item(14, 5).Value = '  +0.30%  '
$ws.Cells.Item(15, 4).Value = '0.0₅5746'
$ws.Cells.Item(15, 5).Value = '  -1.77%  '
$ws.Cells.Item(16, 4).Value = "'" + '64.35'
$ws.Cells.Item(16, 5).Value = '  -0.59%  '
$ws.Cells.Item(17, 4).Value = '26.024.35'
$ws.Cells.Item(17, 5).Value = '  -0.04%  '
$ws.Cells.Item(18, 4).Value = "'" + '1.001'
$ws.Cells.Item(18, 5).Value = '  -0.67%  '
$ws.Cells.Item(19, 4).Value = "'" + '195.89'
$ws.Cells.Item(19, 5).Value = '  -1.38%  '
$ws.Cells.Item(20, 4).Value = "'" + '4.434'
$ws.Cells.Item(21, 4).Value = "'" + '9.951'
$ws.Cells.Item(21, 5).Value = '  -0.67%  '
$ws.Cells.Item(22, 4).Value = "'" + '6.115'
$ws.Cells.Item(22, 5).Value = '  +0.80%  '
$ws.Cells.Item(23, 4).Value = "'" + '1.002'
$ws.Cells.Item(23, 5).Value = '  -0.58%  '
$ws.Cells.Item(24, 5).Value = '  +1.45%  '
$ws.Cells.Item(25, 4).Value = "'" + '143.49'
$ws.Cells.Item(25, 5).Value = '  +2.23%  '
$ws.Cells.Item(26, 4).Value = "'" + '0.1242'
$ws.Cells.Item(26, 5).Value = '  +7.97%  '
$ws.Cells.Item(27, 4).Value = "'" + '6.886'
$ws.Cells.Item(27, 5).Value = '  -0.23%  '
$ws.Cells.Item(28, 4).Value = "'" + '15.65'
$ws.Cells.Item(28, 5).Value = '  -0.65%  '
$ws.Cells.Item(29, 5).Value = '  -0.27%  '
$ws.Cells.Item(30, 4).Value = "'" + '0.04889'
$ws.Cells.Item(30, 5).Value = '  -2.62%  '
$ws.Cells.Item(31, 4).Value = "'" + '3.274'
$ws.Cells.Item(32, 4).Value = "'" + '3.215'
$ws.Cells.Item(32, 5).Value = '  +0.32%  '
$ws.Cells.Item(33, 4).Value = "'" + '1.546'
$ws.Cells.Item(33, 5).Value = '  -0.06%  '
$ws.Cells.Item(34, 4).Value = "'" + '2.375'
$ws.Cells.Item(34, 5).Value = '  +0.36%  '
$ws.Cells.Item(35, 4).Value = "'" + '0.9164'
$ws.Cells.Item(35, 5).Value = '  +2.35%  '
$ws.Cells.Item(36, 4).Value = "'" + '2.571'
$ws.Cells.Item(36, 5).Value = '  -0.94%  '
$ws.Cells.Item(37, 4).Value = "'" + '0.5548'
$ws.Cells.Item(37, 5).Value = '  +0.15%  '
$ws.Cells.Item(38, 4).Value = '1.091.77'
$ws.Cells.Item(38, 5).Value = '  -3.83%  '
$ws.Cells.Item(39, 4).Value = "'" + '0.01573'
$ws.Cells.Item(39, 5).Value = '  +0.48%  '
$ws.Cells.Item(40, 5).Value = '  -0.67%  '
$ws.Cells.Item(41, 2).Value = 'FraxShare'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(41, 4).Value = "'" + '5.604'
$ws.Cells.Item(41, 5).Value = '  -1.02%  '
$ws.Cells.Item(42, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Cells.Item(42, 4).Value = "'" + '0.8061'
$ws.Cells.Item(42, 5).Value = '  -1.36%  '
$ws.Cells.Item(43, 2).Value = 'Quant'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Cells.Item(43, 4).Value = "'" + '98.85'
$ws.Cells.Item(43, 5).Value = '  -1.18%  '
$ws.Cells.Item(44, 2).Value = 'BabyDogeCoin'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Cells.Item(44, 4).Value = '0.0₈120'
$ws.Cells.Item(44, 5).Value = '  -3.87%  '
$ws.Cells.Item(45, 2).Value = 'RocketPoolETH'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Cells.Item(45, 4).Value = '1.779.49'
$ws.Cells.Item(45, 5).Value = '  -0.36%  '
$ws.Cells.Item(46, 2).Value = 'Mantle'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(46, 4).Value = "'" + '0.4537'
$ws.Cells.Item(46, 5).Value = '  -0.09%  '
$ws.Cells.Item(47, 2).Value = 'Aave'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(47, 4).Value = "'" + '55.36'
$ws.Cells.Item(47, 5).Value = '  +0.21%  '
$ws.Cells.Item(48, 2).Value = 'Frax'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Cells.Item(48, 4).Value = "'" + '1.001'
$ws.Cells.Item(48, 5).Value = '  -0.52%  '
$ws.Cells.Item(49, 2).Value = 'Cronos'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Cells.Item(49, 4).Value = "'" + '0.05218'
$ws.Cells.Item(49, 5).Value = '  +2.45%  '
$ws.Cells.Item(50, 2).Value = 'EnergySwap'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(50, 4).Value = "'" + '7.576'
$ws.Cells.Item(50, 5).Value = '  +2.02%  '
$ws.Cells.Item(51, 2).Value = 'USDD'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Cells.Item(51, 4).Value = "'" + '1.004'
$ws.Cells.Item(51, 5).Value = '  -0.44%  '
